$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("data\output\output_subpreg\00093\4003684_p_", "Pregunta no pudo ser procesada"),
    @("data\output\output_subpreg\00206\4007946_p2", "Pregunta no pudo ser procesada"),
    @("data\output\output_subpreg\00206\4007946_p3", "Pregunta no pudo ser procesada"),
    @("data\output\output_subpreg\00206\4007946_p26", "Pregunta no pudo ser procesada"),
    @("data\output\output_subpreg\00206\4007946_p25", "Pregunta no pudo ser procesada")
)

$startRow = 30
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}
